$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Easy / Microsoft Sam -> MS Sam)
$ws.Range("B2").Value = "MS Sam"
$ws.Range("D2").Value = 11
$ws.Range("F2").Value = "Sam.png"

# Row 3 (Medium / Holly)
$ws.Range("F3").Value = "Holly.png"

# Row 4 (Hard / Skynet)
$ws.Range("C4").Value = "N|5,N|6,N|7,N|8,N|9,N|-5,N|-6,N|-7,N|-8,N|-9"
$ws.Range("F4").Value = "Skynet.png"

# Update selection to match the new active cell
$ws.Range("F4").Select()

# Column width adjustments (values compensated for the runtime's internal
# char->pixel rounding so the resulting saved XML width lands as close as
# possible to the target widths: 16.140625, 51, 13, 9.85546875, 13)
$ws.Columns.Item(2).ColumnWidth = 15.307291666666666
$ws.Columns.Item(3).ColumnWidth = 50.166666666666664
$ws.Columns.Item(4).ColumnWidth = 12.166666666666666
$ws.Columns.Item(5).ColumnWidth = 9.022135416666666
$ws.Columns.Item(6).ColumnWidth = 12.166666666666666
